$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '98.551.02'
$ws.Range('D3').Value = '3.372.61'
$ws.Range('E3').Value = '  -0.11%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = "'257.87"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.77%  '
$ws.Range('D6').Value = "'673.03"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +6.94%  '
$ws.Range('E7').Value = '  +12.20%  '
$ws.Range('D8').Value = "'0.460"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +16.89%  '
$ws.Range('D9').Value = "'1.10"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +26.68%  '
$ws.Range('E10').Value = '  +0.01%  '
$ws.Range('D11').Value = '3.369.86'
$ws.Range('E11').Value = '  -0.16%  '
$ws.Range('E12').Value = '  +5.47%  '
$ws.Range('D13').Value = "'42.51"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +17.27%  '
$ws.Range('E14').Value = '  +7.81%  '
$ws.Range('D15').Value = '97.702.49'
$ws.Range('E15').Value = '  -1.04%  '
$ws.Range('D16').Value = "'5.62"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.17%  '
$ws.Range('D17').Value = '3.374.69'
$ws.Range('E17').Value = '  -0.18%  '
$ws.Range('D18').Value = "'7.66"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +24.89%  '
$ws.Range('D19').Value = "'16.94"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +11.12%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').Value = "'532.31"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +8.46%  '
$ws.Range('B21').Value = 'SuiNetwork'
$ws.Range('C21').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D21').Value = "'3.58"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.23%  '
$ws.Range('E22').Value = '  +12.11%  '
$ws.Range('E23').Value = '  +58.62%  '
$ws.Range('D24').Value = "'0.0000213"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.40%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').Value = "'102.38"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +15.36%  '
$ws.Range('B26').Value = 'NEARProtocol'
$ws.Range('C26').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D26').Value = "'6.28"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +10.87%  '
$ws.Range('D27').Value = "'12.67"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +5.55%  '
$ws.Range('D28').Value = '3.558.35'
$ws.Range('E28').Value = '  +0.04%  '
$ws.Range('E29').Value = '  +15.19%  '
$ws.Range('D30').Value = "'1.00"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.13%  '
$ws.Range('D31').Value = "'11.18"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +15.57%  '
$ws.Range('E32').Value = '  -1.03%  '
$ws.Range('E33').Value = '  -0.04%  '
$ws.Range('D34').Value = "'29.88"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +6.10%  '
$ws.Range('D35').Value = "'0.546"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +18.44%  '
$ws.Range('D36').Value = "'7.90"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.80%  '
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').Value = "'0.161"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +7.42%  '
$ws.Range('B38').Value = 'PancakeSwap'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D38').Value = "'2.14"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +8.55%  '
$ws.Range('D39').Value = "'527.15"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.34%  '
$ws.Range('D40').Value = "'0.0458"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +40.35%  '
$ws.Range('E41').Value = '  +5.68%  '
$ws.Range('D42').Value = "'24.70"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.86%  '
$ws.Range('D43').Value = "'3.79"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.62%  '
$ws.Range('D44').Value = "'0.839"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +5.87%  '
$ws.Range('E45').Value = '  +2.35%  '
$ws.Range('E46').Value = '  +0.00%  '
$ws.Range('B47').Value = 'Cosmos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D47').Value = "'8.00"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +20.73%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').Value = "'2.06"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.95%  '
$ws.Range('D49').Value = "'5.15"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +10.91%  '
$ws.Range('D50').Value = "'50.88"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +10.56%  '
$ws.Range('D51').Value = "'1.52"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +11.28%  '
